$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95
$ws.Range("A95").Value = 130964535
$ws.Range("Q95").Value = 509939
$ws.Range("R95").Value = 6719007
$ws.Range("AC95").Value = "Flera . inventering åt vasa vind"

# Row 96
$ws.Range("A96").Value = 130964529
$ws.Range("Q96").Value = 509588
$ws.Range("R96").Value = 6719099
$ws.Range("AC96").Value = "Enstaka . inventering åt vasa vind"
$ws.Range("AX96").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 97
$ws.Range("A97").Value = 130964573
$ws.Range("Q97").Value = 509515
$ws.Range("R97").Value = 6719063
$ws.Range("AC97").Value = "Måttliga förekomster . inventering åt vasa vind"
$ws.Range("AX97").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 99
$ws.Range("AX99").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 100
$ws.Range("AX100").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 101
$ws.Range("AX101").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 102
$ws.Range("AX102").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 103
$ws.Range("AX103").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 104
$ws.Range("A104").Value = 130964642
$ws.Range("B104").Value = 99037
$ws.Range("D104").Value = "LC"
$ws.Range("E104").Value = 221952
$ws.Range("F104").Value = "Spindelblomster"
$ws.Range("G104").Value = "Neottia cordata"
$ws.Range("H104").Value = "(L.) Rich."
$ws.Range("Q104").Value = 509917
$ws.Range("R104").Value = 6719042
$ws.Range("AC104").Value = "Måttliga förekomster . inventering åt vasa vind"
$ws.Range("AX104").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 105
$ws.Range("AX105").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 106
$ws.Range("A106").Value = 130964541
$ws.Range("B106").Value = 91809
$ws.Range("D106").Value = "NT"
$ws.Range("E106").Value = 1202
$ws.Range("F106").Value = "Ullticka"
$ws.Range("G106").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H106").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q106").Value = 509703
$ws.Range("R106").Value = 6719018
$ws.Range("AC106").Value = "Enstaka . inventering åt vasa vind"
$ws.Range("AX106").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 107
$ws.Range("A107").Value = 130964545
$ws.Range("B107").Value = 57073
$ws.Range("E107").Value = 100138
$ws.Range("F107").Value = "Tjäder"
$ws.Range("G107").Value = "Tetrao urogallus"
$ws.Range("H107").Value = "Linnaeus, 1758"
$ws.Range("Q107").Value = 509535
$ws.Range("R107").Value = 6718925
$ws.Range("AC107").Value = "Spillning . inventering åt vasa vind"

# Row 108
$ws.Range("AX108").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 112
$ws.Range("A112").Value = 130964533
$ws.Range("B112").Value = 79244
$ws.Range("D112").Value = "NT"
$ws.Range("E112").Value = 6425
$ws.Range("F112").Value = "Garnlav"
$ws.Range("G112").Value = "Alectoria sarmentosa"
$ws.Range("H112").Value = "(Ach.) Ach."
$ws.Range("Q112").Value = 509984
$ws.Range("R112").Value = 6719028
$ws.Range("AC112").Value = "Rikligt . inventering åt vasa vind"

# Row 113
$ws.Range("A113").Value = 130964645
$ws.Range("B113").Value = 99037
$ws.Range("D113").Value = "LC"
$ws.Range("E113").Value = 221952
$ws.Range("F113").Value = "Spindelblomster"
$ws.Range("G113").Value = "Neottia cordata"
$ws.Range("H113").Value = "(L.) Rich."
$ws.Range("Q113").Value = 509804
$ws.Range("R113").Value = 6719024
$ws.Range("AC113").Value = "Måttliga förekomster . inventering åt vasa vind"
$ws.Range("AX113").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 115
$ws.Range("AX115").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 116
$ws.Range("AX116").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 117
$ws.Range("AX117").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 119
$ws.Range("A119").Value = 130964648
$ws.Range("B119").Value = 92268
$ws.Range("D119").Value = "VU"
$ws.Range("E119").Value = 1209
$ws.Range("F119").Value = "Rynkskinn"
$ws.Range("G119").Value = "Hermanssonia centrifuga"
$ws.Range("H119").Value = "(P. Karst.) Zmitr."
$ws.Range("Q119").Value = 509744
$ws.Range("R119").Value = 6718982
$ws.Range("AC119").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 120
$ws.Range("A120").Value = 130964649
$ws.Range("B120").Value = 98931
$ws.Range("D120").Value = "LC"
$ws.Range("E120").Value = 219790
$ws.Range("F120").Value = "Fläcknycklar"
$ws.Range("G120").Value = "Dactylorhiza maculata"
$ws.Range("H120").Value = "(L.) Soó"
$ws.Range("Q120").Value = 509705
$ws.Range("R120").Value = 6718923
$ws.Range("AC120").Value = "Måttlig förekomst . inventering åt vasa vind"

# Row 121
$ws.Range("AX121").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 122
$ws.Range("A122").Value = 130964542
$ws.Range("B122").Value = 57073
$ws.Range("E122").Value = 100138
$ws.Range("F122").Value = "Tjäder"
$ws.Range("G122").Value = "Tetrao urogallus"
$ws.Range("H122").Value = "Linnaeus, 1758"
$ws.Range("Q122").Value = 509635
$ws.Range("R122").Value = 6718941
$ws.Range("AC122").Value = "Spillning . inventering åt vasa vind"

# Row 124
$ws.Range("A124").Value = 130964644
$ws.Range("B124").Value = 98918
$ws.Range("E124").Value = 220093
$ws.Range("F124").Value = "Korallrot"
$ws.Range("G124").Value = "Corallorhiza trifida"
$ws.Range("H124").Value = "Châtel."
$ws.Range("Q124").Value = 509801
$ws.Range("R124").Value = 6719017
$ws.Range("AC124").Value = "Sparsamma förekomster . inventering åt vasa vind"
$ws.Range("AX124").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 127
$ws.Range("A127").Value = 130964538
$ws.Range("B127").Value = 79244
$ws.Range("D127").Value = "NT"
$ws.Range("E127").Value = 6425
$ws.Range("F127").Value = "Garnlav"
$ws.Range("G127").Value = "Alectoria sarmentosa"
$ws.Range("H127").Value = "(Ach.) Ach."
$ws.Range("Q127").Value = 509875
$ws.Range("R127").Value = 6719025
$ws.Range("AC127").Value = "Enstaka . inventering åt vasa vind"
$ws.Range("AX127").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 128
$ws.Range("A128").Value = 130964546
$ws.Range("B128").Value = 92504
$ws.Range("D128").Value = "VU"
$ws.Range("E128").Value = 898
$ws.Range("F128").Value = "Blackticka"
$ws.Range("G128").Value = "Steccherinum collabens"
$ws.Range("H128").Value = "(Fr.) Vesterholt"
$ws.Range("Q128").Value = 509515
$ws.Range("R128").Value = 6718886
$ws.Range("AC128").Value = "Betydande förekomst . inventering åt vasa vind"
$ws.Range("AX128").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"

# Row 129
$ws.Range("A129").Value = 130964639
$ws.Range("B129").Value = 57073
$ws.Range("D129").Value = "LC"
$ws.Range("E129").Value = 100138
$ws.Range("F129").Value = "Tjäder"
$ws.Range("G129").Value = "Tetrao urogallus"
$ws.Range("H129").Value = "Linnaeus, 1758"
$ws.Range("Q129").Value = 509645
$ws.Range("R129").Value = 6719169
$ws.Range("AC129").Value = "Vinterspillning . inventering åt vasa vind"

# Row 130
$ws.Range("A130").Value = 130964646
$ws.Range("B130").Value = 91809
$ws.Range("D130").Value = "NT"
$ws.Range("E130").Value = 1202
$ws.Range("F130").Value = "Ullticka"
$ws.Range("G130").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H130").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q130").Value = 509764
$ws.Range("R130").Value = 6719043
$ws.Range("AC130").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 131
$ws.Range("AX131").Value = "Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning"
